# This script updates the cryptocurrency price (column D) and 1-hour volume
# change percentage (column E) values on Sheet1, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "30.648.34", "0.9999") that must
# remain text and not be auto-converted to numbers by Excel. We briefly mark each
# target cell as Text before writing the value, then restore its original "Normal"
# style so no visible formatting change is introduced.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '30.648.34'
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.883.92'
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '249.65'
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2943'
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06543'
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '22.05'
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07743'
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '97.04'
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.7396'
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '1.881.24'
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.249'
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '275.35'
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '30.748.70'
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '13.20'
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000007556'
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '2.128.72'
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.360'
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '9.241'
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '164.07'
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.918'
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.09741'
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.506'
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.299'
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.173'
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.04889'
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7006'
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.719'
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01921'
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.793'
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '6.323'
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '75.50'
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.4263'
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.8423'
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '102.74'
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.397'
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.069'
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '918.21'
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05766'
$c.Style = "Normal"

# Column E holds percentage-change text already including a leading/trailing
# double-space padding and +/- sign, e.g. "  +0.61%  ". These are plain strings
# that Excel keeps as text automatically, so no special handling is required.

$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +1.54%  '
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("E18").Value = '  -2.97%  '
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("E22").Value = '  +2.12%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  +0.83%  '
$ws.Range("E25").Value = '  -1.03%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("E29").Value = '  -1.94%  '
$ws.Range("E30").Value = '  -2.10%  '
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("E33").Value = '  +2.46%  '
$ws.Range("E34").Value = '  +2.20%  '
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("E38").Value = '  +2.91%  '
$ws.Range("E39").Value = '  +2.24%  '
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("E41").Value = '  +6.83%  '
$ws.Range("E42").Value = '  +4.57%  '
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("E51").Value = '  +2.10%  '
